# Update countries & provincias Spain
# Applies the data refresh captured by the diff:
#  - Swap the country labels for rows 144/145 (Uganda now ranks above Georgia)
#  - Swap the country labels for rows 210/211 (Groenlandia now ranks above Islas Malvinas)
#  - Refresh the numeric statistics for the affected rows
#  - Update the "Datos actualizados" timestamp

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Swap country names (rank changed between updates) ---
$ws.Range("A144").Value = "Uganda"
$ws.Range("A145").Value = "Georgia"

$ws.Range("A210").Value = "Groenlandia"
$ws.Range("A211").Value = "Islas Malvinas"

# --- Row 6: India ---
$ws.Range("B6").Value = 1195674
$ws.Range("C6").Value = 1589
$ws.Range("E6").Value = 413843
$ws.Range("G6").Value = 11
$ws.Range("H6").Value = 28781

# --- Row 14: Iran ---
$ws.Range("B14").Value = 281413
$ws.Range("C14").Value = 2586
$ws.Range("D14").Value = 244840
$ws.Range("E14").Value = 21720
$ws.Range("G14").Value = 219
$ws.Range("H14").Value = 14853

# --- Row 49: Rumania ---
$ws.Range("B49").Value = 40163
$ws.Range("C49").Value = 1030
$ws.Range("D49").Value = 24663
$ws.Range("E49").Value = 13399
$ws.Range("G49").Value = 27
$ws.Range("H49").Value = 2101

# --- Row 55: Suiza ---
$ws.Range("B55").Value = 33883
$ws.Range("C55").Value = 141
$ws.Range("E55").Value = 1611

# --- Row 66: Nepal ---
$ws.Range("B66").Value = 18094
$ws.Range("C66").Value = 100
$ws.Range("D66").Value = 12684
$ws.Range("E66").Value = 5368
$ws.Range("G66").Value = 2
$ws.Range("H66").Value = 42

# --- Row 87: Consejo Danes para los Refugiados ---
$ws.Range("B87").Value = 8626
$ws.Range("C87").Value = 92
$ws.Range("D87").Value = 4790
$ws.Range("E87").Value = 3639
$ws.Range("G87").Value = 1
$ws.Range("H87").Value = 197

# --- Row 124: Eslovenia ---
$ws.Range("D124").Value = 1648
$ws.Range("E124").Value = 243

# --- Row 144: now Uganda ---
$ws.Range("B144").Value = 1075
$ws.Range("C144").Value = 0
$ws.Range("D144").Value = 958
$ws.Range("E144").Value = 117
$ws.Range("H144").Value = 0

# --- Row 145: now Georgia ---
$ws.Range("B145").Value = 1073
$ws.Range("C145").Value = 24
$ws.Range("D145").Value = 907
$ws.Range("E145").Value = 150
$ws.Range("H145").Value = 16

# --- Row 146: Burkina Faso ---
$ws.Range("B146").Value = 1066
$ws.Range("C146").Value = 1
$ws.Range("D146").Value = 917
$ws.Range("E146").Value = 96

# --- Update timestamp ---
$ws.Range("A1").Value = "Datos actualizados a 22 de Julio de 2020 a las 12:59"
